# Append a new row (row 11, "2021年") to Sheet1, mirroring the format of
# the preceding year row (row 10) and filling in the reported indicator
# values for 2021, leaving columns E and V blank (as in the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 10's formatting down to row 11 first (so the new row picks up
# the same styles - bold/border/centered label cell in column A, plain
# numeric cells elsewhere), then overwrite the values for 2021.
$ws.Range("A10:AQ10").Copy()
$ws.Range("A11:AQ11").PasteSpecial(-4122)

$ws.Range("A11").Value = "2021年"

$ws.Range("B11").Value = 20.54
$ws.Range("C11").Value = 5.45
$ws.Range("D11").Value = 3.05
$ws.Range("F11").Value = 41.1
$ws.Range("G11").Value = 89.2
$ws.Range("H11").Value = 21.98
$ws.Range("I11").Value = 15.58
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 896.87
$ws.Range("L11").Value = 6.28
$ws.Range("M11").Value = 1.87
$ws.Range("N11").Value = 0.51
$ws.Range("O11").Value = 15.01
$ws.Range("P11").Value = 74.31999999999999
$ws.Range("Q11").Value = 3.8
$ws.Range("R11").Value = 3.09
$ws.Range("S11").Value = 28.02
$ws.Range("T11").Value = 27.35
$ws.Range("U11").Value = -89.64
$ws.Range("W11").Value = 31.58
$ws.Range("X11").Value = 9.029999999999999
$ws.Range("Y11").Value = 211.78
$ws.Range("Z11").Value = 60.52
$ws.Range("AA11").Value = 12.62
$ws.Range("AB11").Value = 18.11
$ws.Range("AC11").Value = 20.95
$ws.Range("AD11").Value = 30.3
$ws.Range("AE11").Value = 22.93
$ws.Range("AF11").Value = 41.05
$ws.Range("AG11").Value = 17.13
$ws.Range("AH11").Value = 49.12
$ws.Range("AI11").Value = -3.03
$ws.Range("AJ11").Value = 4.47
$ws.Range("AK11").Value = 30.21
$ws.Range("AL11").Value = 7.4
$ws.Range("AM11").Value = 22.78
$ws.Range("AN11").Value = 0.64
$ws.Range("AO11").Value = -8.800000000000001
$ws.Range("AP11").Value = 46.12
$ws.Range("AQ11").Value = 1.96

# E11 and V11 stay blank (empty string), matching the source row.
$ws.Range("E11").Value = ""
$ws.Range("V11").Value = ""
